$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Number" -> "Composition" in G1
$ws.Range("G1").Value = "Composition"

# Apply right alignment to the whole column G (matches the new column style)
$ws.Range("G1:G1048576").HorizontalAlignment = -4152

# Select the column to mirror the saved selection state
$ws.Range("G1:G1048576").Select()
